$d = $word.ActiveDocument

# Update the date line (first paragraph)
$d.Paragraphs.Item(1).Range.Text = "2024-02-22 Thursday"

# Update the 20x5 table of addition/subtraction facts
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "22+43=65"
$t.Cell(1, 2).Range.Text = "33-9=24"
$t.Cell(1, 3).Range.Text = "79-57=22"
$t.Cell(1, 4).Range.Text = "10+41=51"
$t.Cell(1, 5).Range.Text = "2+73=75"

$t.Cell(2, 1).Range.Text = "39+23=62"
$t.Cell(2, 2).Range.Text = "5+79=84"
$t.Cell(2, 3).Range.Text = "46-22=24"
$t.Cell(2, 4).Range.Text = "66-22=44"
$t.Cell(2, 5).Range.Text = "96-91=5"

$t.Cell(3, 1).Range.Text = "91-47=44"
$t.Cell(3, 2).Range.Text = "9+90=99"
$t.Cell(3, 3).Range.Text = "96-7=89"
$t.Cell(3, 4).Range.Text = "9+77=86"
$t.Cell(3, 5).Range.Text = "45+7=52"

$t.Cell(4, 1).Range.Text = "42-3=39"
$t.Cell(4, 2).Range.Text = "66-7=59"
$t.Cell(4, 3).Range.Text = "75+22=97"
$t.Cell(4, 4).Range.Text = "1+45=46"
$t.Cell(4, 5).Range.Text = "65-23=42"

$t.Cell(5, 1).Range.Text = "27+65=92"
$t.Cell(5, 2).Range.Text = "79-44=35"
$t.Cell(5, 3).Range.Text = "95+0=95"
$t.Cell(5, 4).Range.Text = "70-55=15"
$t.Cell(5, 5).Range.Text = "67+20=87"

$t.Cell(6, 1).Range.Text = "31+53=84"
$t.Cell(6, 2).Range.Text = "42-29=13"
$t.Cell(6, 3).Range.Text = "93-9=84"
$t.Cell(6, 4).Range.Text = "10+9=19"
$t.Cell(6, 5).Range.Text = "39+60=99"

$t.Cell(7, 1).Range.Text = "1+56=57"
$t.Cell(7, 2).Range.Text = "59+9=68"
$t.Cell(7, 3).Range.Text = "1+39=40"
$t.Cell(7, 4).Range.Text = "80-19=61"
$t.Cell(7, 5).Range.Text = "74-51=23"

$t.Cell(8, 1).Range.Text = "10+10=20"
$t.Cell(8, 2).Range.Text = "91-35=56"
$t.Cell(8, 3).Range.Text = "44+55=99"
$t.Cell(8, 4).Range.Text = "75-63=12"
$t.Cell(8, 5).Range.Text = "26+61=87"

$t.Cell(9, 1).Range.Text = "51+19=70"
$t.Cell(9, 2).Range.Text = "35+36=71"
$t.Cell(9, 3).Range.Text = "92-57=35"
$t.Cell(9, 4).Range.Text = "45+4=49"
$t.Cell(9, 5).Range.Text = "38+2=40"

$t.Cell(10, 1).Range.Text = "67-9=58"
$t.Cell(10, 2).Range.Text = "55-9=46"
$t.Cell(10, 3).Range.Text = "4+17=21"
$t.Cell(10, 4).Range.Text = "32+53=85"
$t.Cell(10, 5).Range.Text = "45-29=16"

$t.Cell(11, 1).Range.Text = "59+16=75"
$t.Cell(11, 2).Range.Text = "87-75=12"
$t.Cell(11, 3).Range.Text = "88-19=69"
$t.Cell(11, 4).Range.Text = "64+10=74"
$t.Cell(11, 5).Range.Text = "57-54=3"

$t.Cell(12, 1).Range.Text = "37+47=84"
$t.Cell(12, 2).Range.Text = "63-53=10"
$t.Cell(12, 3).Range.Text = "53-41=12"
$t.Cell(12, 4).Range.Text = "8+30=38"
$t.Cell(12, 5).Range.Text = "45+21=66"

$t.Cell(13, 1).Range.Text = "98-66=32"
$t.Cell(13, 2).Range.Text = "23+42=65"
$t.Cell(13, 3).Range.Text = "23+24=47"
$t.Cell(13, 4).Range.Text = "71-16=55"
$t.Cell(13, 5).Range.Text = "6+22=28"

$t.Cell(14, 1).Range.Text = "35-21=14"
$t.Cell(14, 2).Range.Text = "80-4=76"
$t.Cell(14, 3).Range.Text = "10+66=76"
$t.Cell(14, 4).Range.Text = "38-20=18"
$t.Cell(14, 5).Range.Text = "69-21=48"

$t.Cell(15, 1).Range.Text = "53-47=6"
$t.Cell(15, 2).Range.Text = "22+3=25"
$t.Cell(15, 3).Range.Text = "42+43=85"
$t.Cell(15, 4).Range.Text = "2+40=42"
$t.Cell(15, 5).Range.Text = "41+23=64"

$t.Cell(16, 1).Range.Text = "85-59=26"
$t.Cell(16, 2).Range.Text = "8+91=99"
$t.Cell(16, 3).Range.Text = "53-29=24"
$t.Cell(16, 4).Range.Text = "81-12=69"
$t.Cell(16, 5).Range.Text = "92-21=71"

$t.Cell(17, 1).Range.Text = "0+35=35"
$t.Cell(17, 2).Range.Text = "73-72=1"
$t.Cell(17, 3).Range.Text = "19+50=69"
$t.Cell(17, 4).Range.Text = "92-78=14"
$t.Cell(17, 5).Range.Text = "52-46=6"

$t.Cell(18, 1).Range.Text = "29+51=80"
$t.Cell(18, 2).Range.Text = "16+24=40"
$t.Cell(18, 3).Range.Text = "5+16=21"
$t.Cell(18, 4).Range.Text = "76-45=31"
$t.Cell(18, 5).Range.Text = "65-8=57"

$t.Cell(19, 1).Range.Text = "43+9=52"
$t.Cell(19, 2).Range.Text = "3+74=77"
$t.Cell(19, 3).Range.Text = "42+14=56"
$t.Cell(19, 4).Range.Text = "97-92=5"
$t.Cell(19, 5).Range.Text = "97-18=79"

$t.Cell(20, 1).Range.Text = "12-1=11"
$t.Cell(20, 2).Range.Text = "37-19=18"
$t.Cell(20, 3).Range.Text = "45+7=52"
$t.Cell(20, 4).Range.Text = "66+22=88"
$t.Cell(20, 5).Range.Text = "22+44=66"
